# Apply "cap nhat nang suat" updates to the timesheet workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (row 2)
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0

# Row 3
$ws.Range("F3").Value = 0

# Row 7
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0

# Row 8
$ws.Range("F8").Value = 0

# Row 14
$ws.Range("G14").Value = 0

# Row 16
$ws.Range("F16").Value = 1

# Row 17
$ws.Range("F17").Value = 0.5

# Row 19
$ws.Range("F19").Value = 0

# Row 28
$ws.Range("G28").Value = 0

# Update the frozen-pane view and active selection on the sheet view
$ws.Activate()
$ws.Range("D33").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 8
$win.ScrollColumn = 2
